# Upgraded AttendanceBelow75Activity for college
#
# This script applies the following changes to the document:
#  1. Remove the paragraph "Change UI to accept `semester` of student"
#  2. Remove the paragraph "Unit test marks"
#  3. Remove the paragraph "Submission -> micro project, manual, oral"
#  4. Change "Semester 1 engineering graphics subject code" to
#     "All branches, how many division and batches"
#  5. Remove the paragraphs from "CO5-A_count" through "subject_short_name"
#  6. Clear the text "teacher_of_semester_1 : true" from the final paragraph,
#     while keeping the paragraph (and its bookmark) intact.

$d = $word.ActiveDocument

function Get-ParaText($p) {
    # Range.Text includes the trailing paragraph-mark character; strip it
    # off so comparisons against plain strings work as expected.
    return $p.Range.Text.TrimEnd([char]13)
}

function Find-ParagraphByText($text) {
    foreach ($p in $d.Paragraphs) {
        if ((Get-ParaText $p) -eq $text) {
            return $p
        }
    }
    return $null
}

function Find-ParagraphAfter($text) {
    $paras = @()
    foreach ($p in $d.Paragraphs) { $paras += $p }
    for ($i = 0; $i -lt $paras.Count; $i++) {
        if ((Get-ParaText $paras[$i]) -eq $text) {
            return $paras[$i + 1]
        }
    }
    return $null
}

# Delete an entire paragraph (including its paragraph mark) by extending
# the deletion range up to the start of the paragraph that follows it.
function Remove-WholeParagraphRange($startPara, $stopPara) {
    $r = $d.Range($startPara.Range.Start, $stopPara.Range.Start)
    $r.Delete()
}

# --- Work from the bottom of the document upward so paragraph indices /
# --- object references for not-yet-processed paragraphs stay valid. ---

# 6. Final paragraph: remove only the run text "teacher_of_semester_1 : true"
#    but keep the paragraph mark and the _GoBack bookmark that follows it.
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$full = $pLast.Range
$r = $d.Range($full.Start, $full.End - 1)
$r.Delete()

# 5. Remove paragraphs "CO5-A_count" .. "subject_short_name" (8 paragraphs),
#    i.e. everything from "CO5-A_count" up to (but not including) the final
#    (now-empty) paragraph that used to hold "teacher_of_semester_1 : true".
$pCOStart = Find-ParagraphByText "CO5-A_count"
$pFinal = $d.Paragraphs($d.Paragraphs.Count)
Remove-WholeParagraphRange $pCOStart $pFinal

# 4. Update text of the "Semester 1 engineering graphics subject code" line.
$d.Content.Find.Execute("Semester 1 engineering graphics subject code", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "All branches, how many division and batches", 2)

# 3 & 2. Remove "Unit test marks" and "Submission -> micro project, manual, oral"
#    paragraphs (they are adjacent, so remove as a single block).
$pUnitTest = Find-ParagraphByText "Unit test marks"
$pAfterSubmission = Find-ParagraphAfter "Submission -> micro project, manual, oral"
Remove-WholeParagraphRange $pUnitTest $pAfterSubmission

# 1. Remove "Change UI to accept `semester` of student" paragraph.
$pChangeUI = Find-ParagraphByText "Change UI to accept ``semester`` of student"
$pAfterChangeUI = Find-ParagraphAfter "Change UI to accept ``semester`` of student"
Remove-WholeParagraphRange $pChangeUI $pAfterChangeUI
